$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.106.31'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.446.07'
$ws.Range('E3').Value = '  +3.45%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '575.30'
$ws.Range('E5').Value = '  +3.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.06'
$ws.Range('E6').Value = '  +3.64%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.456.22'
$ws.Range('E8').Value = '  +3.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.550'
$ws.Range('E9').Value = '  +3.81%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.53'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.124'
$ws.Range('E11').Value = '  +5.04%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.442'
$ws.Range('E12').Value = '  +1.67%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.060.50'
$ws.Range('E13').Value = '  +3.89%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000194'
$ws.Range('E15').Value = '  +8.08%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.47'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.224.72'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.448.03'
$ws.Range('E18').Value = '  +3.31%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.40'
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.32'
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '392.53'
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.42'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '73.03'
$ws.Range('E23').Value = '  +3.46%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.543'
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000121'
$ws.Range('E26').Value = '  +26.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.69'
$ws.Range('E27').Value = '  +10.39%  '
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.11'
$ws.Range('E30').Value = '  +9.87%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.03'
$ws.Range('E31').Value = '  +2.84%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.37'
$ws.Range('E32').Value = '  +6.46%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.68'
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.55'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.03'
$ws.Range('E36').Value = '  +5.55%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '160.68'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.46'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0784'
$ws.Range('E39').Value = '  +6.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '27.50'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.934.66'
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.85'
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0320'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.44'
$ws.Range('E44').Value = '  +2.99%  '
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '42.02'
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '23.77'
$ws.Range('E47').Value = '  +8.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.08'
$ws.Range('E48').Value = '  +4.86%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.20'
$ws.Range('E49').Value = '  +24.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.860'
$ws.Range('E50').Value = '  +6.95%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.55'
$ws.Range('E51').Value = '  +4.54%  '
